$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2 updates ---
$ws.Range("A2").Value = "export"
$ws.Range("C2").Value = "test"
$ws.Range("D2").Value = "'2"
$ws.Range("E2").Value = 44594

# --- Row 3 updates ---
$ws.Range("C3").Value = "test2"
$ws.Range("D3").Value = "'1"
$ws.Range("E3").Value = 44844

# --- Remove row 4 entirely (data + formatting), which shrinks the
#     sheet dimension to A1:E3 and auto-shrinks the A2:A4 merge to A2:A3 ---
$ws.Rows(4).Delete()
